# "update executable build flow" — synchronize main stream with MKT/Tech
# build flow (ABLA core).
#
# 1) The "C Compiler" rounded-rectangle box (inside the top-level "Group 6")
#    gains a second run " (WASI SDK)" so the label reads "C Compiler (WASI SDK)".
# 2) The now-obsolete "Managed C / last version 2023 11 16" caption textbox
#    (TextBox 32) is removed from the slide.
# 3) Best-effort: the handout/notes master "Date Placeholder" fields were
#    re-cached by PowerPoint from "décembre 24" to "novembre 25" on save;
#    attempt the same update (wrapped defensively since master placeholder
#    fields are not always reachable for edits on every host).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) "C Compiler" -> "C Compiler (WASI SDK)" -----------------------------
$topGroup = $s.Shapes.Item(1)
$cCompilerShape = $topGroup.GroupItems.Item("Rounded Rectangle 22")
[void]$cCompilerShape.TextFrame.TextRange.InsertAfter(" (WASI SDK)")

# --- 2) Delete the stale "Managed C / last version ..." textbox -------------
$s.Shapes.Item("TextBox 32").Delete()

# --- 3) Refresh the cached handout/notes master date fields (best effort) ---
try {
    $handoutMaster = $p.HandoutMaster
    $handoutDate = $handoutMaster.Shapes.Item("Date Placeholder 2")
    $handoutDate.TextFrame.TextRange.Text = "novembre 25"
} catch {
    Write-Output ("handout master date field not updatable: " + $_)
}

try {
    $notesMaster = $p.NotesMaster
    $notesDate = $notesMaster.Shapes.Item("Date Placeholder 2")
    $notesDate.TextFrame.TextRange.Text = "novembre 25"
} catch {
    Write-Output ("notes master date field not updatable: " + $_)
}
